$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.843.88'
$ws.Range("E2").Value = '  +4.68%  '

$ws.Range("D3").Value = '2.631.13'
$ws.Range("E3").Value = '  +5.08%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.16'
$ws.Range("E5").Value = '  +1.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.73'
$ws.Range("E6").Value = '  +1.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  +1.84%  '

$ws.Range("D9").Value = '2.629.56'
$ws.Range("E9").Value = '  +5.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  +7.25%  '

$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("E12").Value = '  +3.21%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.07'
$ws.Range("E13").Value = '  +1.44%  '

$ws.Range("E14").Value = '  +5.85%  '

$ws.Range("E15").Value = '  +6.40%  '

$ws.Range("D16").Value = '72.627.69'
$ws.Range("E16").Value = '  +4.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.69'
$ws.Range("E17").Value = '  +3.16%  '

$ws.Range("D18").Value = '2.630.11'
$ws.Range("E18").Value = '  +4.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '384.58'
$ws.Range("E19").Value = '  +6.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.04'
$ws.Range("E20").Value = '  +6.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.56'
$ws.Range("E21").Value = '  +5.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.18'
$ws.Range("E22").Value = '  +2.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.01'
$ws.Range("E23").Value = '  +19.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.14'
$ws.Range("E24").Value = '  +3.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.39'
$ws.Range("E26").Value = '  +3.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.88'
$ws.Range("E27").Value = '  +9.60%  '

$ws.Range("D28").Value = '2.762.09'
$ws.Range("E28").Value = '  +4.68%  '

$ws.Range("E29").Value = '  +0.18%  '

$ws.Range("E30").Value = '  +7.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '530.18'
$ws.Range("E31").Value = '  +3.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.05'
$ws.Range("E32").Value = '  +4.17%  '

$ws.Range("E33").Value = '  +7.80%  '

$ws.Range("E34").Value = '  +2.98%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.27'
$ws.Range("E36").Value = '  +0.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.35'
$ws.Range("E37").Value = '  +3.23%  '

$ws.Range("E38").Value = '  +1.11%  '

$ws.Range("E39").Value = '  +7.17%  '

$ws.Range("E40").Value = '  -7.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.83'
$ws.Range("E41").Value = '  +6.14%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.07'
$ws.Range("E42").Value = '  +5.50%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  +12.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.333'
$ws.Range("E45").Value = '  +4.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.52'
$ws.Range("E46").Value = '  +1.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.28'
$ws.Range("E47").Value = '  +1.26%  '

$ws.Range("E48").Value = '  +3.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.543'
$ws.Range("E49").Value = '  +5.67%  '

$ws.Range("E50").Value = '  +8.62%  '

$ws.Range("D51").Value = '0.0₆0264'
$ws.Range("E51").Value = '  +5.04%  '
